$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Update employee NIK and name (row 2: C2/D2)
$ws.Range("C2").Value = "EN-4-025"
$ws.Range("D2").Value = "Ade Nurjaya"

# Move the active selection to E4
$ws.Range("E4").Select()
